$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.035.94'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '3.122.43'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.38%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.114.06'
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("E10").Value = '  -3.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = '  -3.34%  '
$ws.Range("E13").Value = '  -4.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.33%  '
$ws.Range("D15").Value = '3.634.99'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("E16").Value = '  +2.39%  '
$ws.Range("D17").Value = '63.048.05'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").Value = '3.118.86'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '472.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.32%  '
$ws.Range("E22").Value = '  -3.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.37%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.58%  '
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.108'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.89%  '
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("E35").Value = '  -2.61%  '
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '0.0₃0704'
$ws.Range("E38").Value = '  -7.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '421.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.64%  '
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("E41").Value = '  -1.24%  '
$ws.Range("E42").Value = '  -11.43%  '
$ws.Range("D43").Value = '2.901.55'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  -5.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.262'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.52%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.112'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.34'
$ws.Range("D51").Style = "Normal"
